$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for the last column (E1): "bóng rổ" -> "bóng chuyền da"
$ws.Range("E1").Value = "bóng chuyền da"

# Update score values in column E (rows 2-70)
$ws.Range("E2").Value = 9
$ws.Range("E3").Value = 7
$ws.Range("E4").Value = 9
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 9
$ws.Range("E7").Value = 10
$ws.Range("E8").Value = 8
$ws.Range("E9").Value = 9
$ws.Range("E10").Value = 9
$ws.Range("E11").Value = 9
$ws.Range("E12").Value = 9
$ws.Range("E13").Value = 9
$ws.Range("E14").Value = 9
$ws.Range("E15").Value = 9
$ws.Range("E16").Value = 9
$ws.Range("E17").Value = 9
$ws.Range("E18").Value = 9
$ws.Range("E19").Value = 9
$ws.Range("E20").Value = 9
$ws.Range("E21").Value = 9
$ws.Range("E22").Value = 9
$ws.Range("E23").Value = 9
$ws.Range("E24").Value = 9
$ws.Range("E25").Value = 9
$ws.Range("E26").Value = 9
$ws.Range("E27").Value = 9
$ws.Range("E28").Value = 9
$ws.Range("E29").Value = 9
$ws.Range("E30").Value = 9
$ws.Range("E31").Value = 9
$ws.Range("E32").Value = 9
$ws.Range("E33").Value = 9
$ws.Range("E34").Value = 9
$ws.Range("E35").Value = 9
$ws.Range("E36").Value = 9
$ws.Range("E37").Value = 9
$ws.Range("E38").Value = 9
$ws.Range("E39").Value = 9
$ws.Range("E40").Value = 9
$ws.Range("E41").Value = 9
$ws.Range("E42").Value = 9
$ws.Range("E43").Value = 9
$ws.Range("E44").Value = 9
$ws.Range("E45").Value = 9
$ws.Range("E46").Value = 9
$ws.Range("E47").Value = 9
$ws.Range("E48").Value = 9
$ws.Range("E49").Value = 9
$ws.Range("E50").Value = 9
$ws.Range("E51").Value = 9
$ws.Range("E52").Value = 9
$ws.Range("E53").Value = 9
$ws.Range("E54").Value = 9
$ws.Range("E55").Value = 9
$ws.Range("E56").Value = 9
$ws.Range("E57").Value = 9
$ws.Range("E58").Value = 9
$ws.Range("E59").Value = 9
$ws.Range("E60").Value = 9
$ws.Range("E61").Value = 9
$ws.Range("E62").Value = 9
$ws.Range("E63").Value = 9
$ws.Range("E64").Value = 9
$ws.Range("E65").Value = 9
$ws.Range("E66").Value = 9
$ws.Range("E67").Value = 9
$ws.Range("E68").Value = 9
$ws.Range("E69").Value = 9
$ws.Range("E70").Value = 9

# Move the active selection to E1 (was E54)
$ws.Range("E1").Select()
